$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update rows 14-58: gebiedscode (B) and naam_kort/naam (C/D)
$data = @(
    @{Row=14; B='elz14'; C='Haspengouw'}
    @{Row=15; B='elz15'; C='Herkenrode'}
    @{Row=16; B='elz16'; C='Houtland en Polder'}
    @{Row=17; B='elz12'; C='Kemp en Duin'}
    @{Row=18; B='elz41'; C='Kempenland'}
    @{Row=19; B='elz17'; C='Klein-Brabant Vaartland'}
    @{Row=20; B='elz19'; C='Leuven'}
    @{Row=21; B='elz20'; C='Leuven Noord'}
    @{Row=22; B='elz21'; C='Leuven Zuid'}
    @{Row=23; B='elz22'; C='Maasland'}
    @{Row=24; B='elz23'; C='Mechelen-Katelijne'}
    @{Row=25; B='elz39'; C='Midden WVL'}
    @{Row=26; B='elz25'; C='Middenkempen'}
    @{Row=27; B='elz26'; C='MidWestLim'}
    @{Row=28; B='elz29'; C='N-O-Waasland'}
    @{Row=29; B='elz28'; C='Noord-Limburg'}
    @{Row=30; B='elz27'; C='Noorderkempen'}
    @{Row=31; B='elz32'; C='Oost-Meetjesland'}
    @{Row=32; B='elz30'; C='Oostende-Bredene'}
    @{Row=33; B='elz31'; C='Oostkust'}
    @{Row=34; B='elz33'; C='Pajottenland'}
    @{Row=35; B='elz34'; C='Pallieterland'}
    @{Row=36; B='elz35'; C='Panacea'}
    @{Row=37; B='elz36'; C='Regio Aalst'}
    @{Row=38; B='elz37'; C='Regio Grimbergen'}
    @{Row=39; B='elz18'; C='Regio Kortrijk'}
    @{Row=40; B='elz24'; C='Regio Menen'}
    @{Row=41; B='elz50'; C='Regio Waregem'}
    @{Row=42; B='elz40'; C='RITS'}
    @{Row=43; B='elz43'; C='RupeLaar'}
    @{Row=44; B='elz44'; C='Schelde en Leie'}
    @{Row=45; B='elz45'; C='Scheldekracht'}
    @{Row=46; B='elz48'; C='Vlaamse Ardennen'}
    @{Row=47; B='elz49'; C='Voorkempen'}
    @{Row=48; B='elz4'; C='WE40'}
    @{Row=49; B='elz51'; C='West-Limburg'}
    @{Row=50; B='elz52'; C='West-Meetjesland'}
    @{Row=51; B='elz38'; C='Westhoek'}
    @{Row=52; B='elz47'; C='Westkust&Polder'}
    @{Row=53; B='elz57'; C='Z-W-Waasland'}
    @{Row=54; B='elz53'; C='Zennevallei'}
    @{Row=55; B='elz54'; C='ZOLim'}
    @{Row=56; B='elz55'; C='ZORA'}
    @{Row=57; B='elz56'; C='Zuiderkempen'}
    @{Row=58; B='elz46'; C='Zuidoost Hageland'}
)

foreach ($item in $data) {
    $ws.Cells.Item($item.Row, 2).Value = $item.B
    $ws.Cells.Item($item.Row, 3).Value = $item.C
    $ws.Cells.Item($item.Row, 4).Value = $item.C
}

# Update rows 59-61: rename "Eerstelijnszone onbekend..." to "ELZ onbekend..."
$ws.Cells.Item(59, 3).Value = 'ELZ onbekend'
$ws.Cells.Item(59, 4).Value = 'ELZ onbekend'
$ws.Cells.Item(60, 3).Value = 'ELZ onbekend (Brussel)'
$ws.Cells.Item(60, 4).Value = 'ELZ onbekend (Brussel)'
$ws.Cells.Item(61, 3).Value = 'ELZ onbekend (Vlaanderen)'
$ws.Cells.Item(61, 4).Value = 'ELZ onbekend (Vlaanderen)'
